$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cards")

# --- Add 4 new cards into the previously-empty rows 27-30 ---

# Row 27: Negotiator
$ws.Range("E27").Value = "Negotiator"
$ws.Range("F27").Value = "MINION"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = "While this unit is on the field - you do not pay wages for your units."

# Row 28: Toxic Frog
$ws.Range("E28").Value = "Toxic Frog"
$ws.Range("F28").Value = "MINION"
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 2
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = "You can return this card to your hand, place 3 poison counters on an enemy unit."

# Row 29: Betrayal
$ws.Range("E29").Value = "Betrayal"
$ws.Range("F29").Value = "UTILITY"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = "Kill a Unit you control, then summon a Unit from your deck with the same cost."

# --- Reword "Burly Zombie" resurrection effect text ---
$ws.Range("J20").Value = "When this minion is killed, you can bury it instead. When this minion is ressurected, its strength and health becomes 7."

# --- Rename existing card "Dissappearing frog" -> "Illusionary frog" ---
$ws.Range("E21").Value = "Illusionary frog"

# Row 30: Bloated Body
$ws.Range("E30").Value = "Bloated Body"
$ws.Range("F30").Value = "MINION"
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = "When this unit is killed, deal 3 damage to all other units on the battlefield."
